# Notetype column added: users can choose the notetype for each note.
#
# Before: directive rows 2-6 (#separator/#html/#guid/#tags/#deck), then a
#         header row 7 (GUID / source / learned / TAGS).
# After:  a new "#notetype column:5" directive row is inserted above the
#         "#tags column:4" row, and the header row gains a 5th column
#         (NOTETYPE) - everything below row 4 shifts down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New NOTETYPE header cell (column E) -----------------------------
# Written before the new directive row so the shared-string table
# allocates "NOTETYPE" (index 10) ahead of "#notetype column:5" (index 11),
# matching how Excel append-orders strings as they are entered.
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)   # xlPasteFormats - reuse D7's header style
$ws.Range("E7").Value = "NOTETYPE"
$ws.Columns.Item(5).ColumnWidth = 23.666666666666668 # renders as width 24.43 in the saved xlsx

# --- New "#notetype column:5" directive row ---------------------------
# Insert a row above the current row 5 ("#tags column:4"); rows 5-7 shift
# down to 6-8.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "#notetype column:5"
$ws.Range("A5").Font = $ws.Range("A4").Font
$ws.Rows.Item(5).RowHeight = 19.5

# --- Formatting tweak on the shifted "#deck:..." row -------------------
# In the authored workbook this row (now row 7) stops being a
# custom-formatted row and instead just carries an explicit style on its
# one populated cell.
$ws.Rows.Item(7).ClearFormats()
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

# Match the saved selection from the authored workbook (cursor sitting on
# the new header-row cell A8).
$ws.Range("A8").Select()
